$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, pushing existing rows 64:77 down to 65:78
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly record
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44995
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = 100112043
$ws.Cells.Item(64, 7).Value = "Pepino dulce"
$ws.Cells.Item(64, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 80
$ws.Cells.Item(64, 11).Value = 20000
$ws.Cells.Item(64, 12).Value = 20000
$ws.Cells.Item(64, 13).Value = 20000
$ws.Cells.Item(64, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 16).Value = 1111
$ws.Cells.Item(64, 17).Value = 18
$ws.Cells.Item(64, 18).Value = "Hortaliza"
